$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing parameter value (spatial_freq values)
$ws.Range("C2").Value = "[0.02222, 0.04444]"

# Add new "pilot" data row (row 3), filling in the order that reproduces
# the expected shared-string append sequence
$ws.Range("H3").Value = "pilot"
$ws.Range("C3").Value = "[0.04]"

# Rename header "orientation" -> "direction"
$ws.Range("A1").Value = "direction"

$ws.Range("B3").Value = "[2]"
$ws.Range("A3").Value = "[180, 147.27, 114.54, 81.81, 49.09, 16.36, 0, -16.36, -49.09, -81.81, -114.54, -147.27]"

$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 6

# Move active selection to A3
[void]$ws.Range("A3").Select()
